# Auto-generated Excel COM-interop edit script
# Applies the cell-value updates described by the commit diff
# (cryptos list refresh, GitHub Actions run on 2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.573.99'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.647.06'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''212.54'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').Value = '''0.536'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.09%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''23.59'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.70%  '
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = '''0.0889'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').Value = '1.880.00'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '1.638.57'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('D15').Value = '''4.04'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').Value = '''64.48'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').Value = '27.538.84'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '''231.18'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.17%  '
$ws.Range('D19').Value = '0.0₃0723'
$ws.Range('E19').Value = '  -0.74%  '
$ws.Range('D20').Value = '''7.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').Value = '''9.74'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.88%  '
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('D25').Value = '''149.01'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('D26').Value = '''7.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '''15.59'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.93%  '
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('E31').Value = '  -3.09%  '
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('D34').Value = '1.424.54'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('D38').Value = '''0.885'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.75%  '
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '''0.817'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').Value = '''65.10'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -6.28%  '
$ws.Range('D46').Value = '1.789.45'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = '''88.18'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.79'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0994'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.63%  '
